$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 41: "Interpretation of LL(k)" slide - change "LL(1) grammar" to
# "LL(k) grammar" in the paragraph describing practical grammars.
# -----------------------------------------------------------------------
$slide41 = $p.Slides.Item(41)
$body41 = $slide41.Shapes.Item(4).TextFrame.TextRange

$para = $body41.Paragraphs(5)
$fullLen = $para.Length
$para.Characters(1, $fullLen).Text = "In practice, the syntax of most programming languages can be defined, or at least closely approximated, by an LL(k) grammar"

# -----------------------------------------------------------------------
# Slide 43: "Recursive Decent Parsing" slide
# -----------------------------------------------------------------------
$slide43 = $p.Slides.Item(43)
$content = $slide43.Shapes.Item(2)

# Resize/reposition the content placeholder slightly (widened box). Top
# and Height are left untouched since they don't change.
$content.Left = 36.12496062992126
$content.Width = 655.2

$tr = $content.TextFrame.TextRange

# --- Paragraph 2 ("parseLoop()  // called when parsing the outer loop") ---
# Merge the "parseLoop" run and the "()... outer loop" run into a single
# run, renaming the function to "parseLoopStmt" and dropping the stray
# spell-check flag that was attached to the old "parseLoop" run. We do
# this by inserting the full replacement word in front of the second run
# (so it picks up that run's clean formatting) and then deleting the
# original flagged "parseLoop" run entirely.
$outer = $tr.Paragraphs(2)
$outerTail = $outer.Characters(10, $outer.Length - 9)   # "()         // called when parsing the outer loop"
[void]$outerTail.InsertBefore("parseLoopStmt")
$outer2 = $tr.Paragraphs(2)
$outerHead = $outer2.Characters(1, 9)                    # original "parseLoop"
$outerHead.Text = ""

# --- Paragraph 6 ("      parseLoop()   // called when paring the inner loop") ---
$inner = $tr.Paragraphs(6)

# Fix the "paring" -> "parsing" typo and, in the same stroke, split this
# phrase out into its own run.
$typo = $inner.Characters(31, 12)
$typo.Text = "when parsing "

$inner2 = $tr.Paragraphs(6)
$innerMid = $inner2.Characters(16, 15)                   # "()   // called "
[void]$innerMid.InsertBefore("Stmt")

$inner3 = $tr.Paragraphs(6)
$innerHead = $inner3.Characters(1, 34)                   # "      parseLoopStmt()   // called "
$innerHead.Text = $innerHead.Text
